$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.036.40'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.93%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.902.35'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.50%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7452'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.77'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.73%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9993'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3074'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.52%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.98'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06904'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08011'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7586'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.908.49'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.248'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.25'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.168'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.045.00'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.07'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.77%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '236.74'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -6.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.165.82'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9990'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9997'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.094'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +6.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.329'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.21'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.84'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1267'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.049'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -6.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.358'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.51%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.528'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.299'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.045'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05358'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.79%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.291'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7425'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.718'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01943'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.762'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.251'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4461'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.77%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.98'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.962'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9990'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8304'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.712'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.31'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.792'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.85%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.063.67'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.66'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1167'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.32%  '
